$wb = $excel.ActiveWorkbook

# --- Rename "Flow Chart template tables" -> "Flow Chart Template Tables" ---
$flowChartTables = $wb.Worksheets.Item("Flow Chart template tables")
$flowChartTables.Name = "Flow Chart Template Tables"

# --- Reorder merged cells on "Proposal (By Plan)" sheet ---
# Move the first 5 merge ranges to the end of the merge list by
# un-merging and re-merging them in the desired order.
$proposalByPlan = $wb.Worksheets.Item("Proposal (By Plan)")
$proposalByPlanCells = @("G7:L7", "N7:T7", "G13:L13", "N13:T13", "C35:T35")
foreach ($c in $proposalByPlanCells) {
    $proposalByPlan.Range($c).UnMerge()
    $proposalByPlan.Range($c).Merge()
}

# --- Reorder merged cells on "Proposal" sheet ---
# Move the first 8 merge ranges to the end of the merge list.
$proposal = $wb.Worksheets.Item("Proposal")
$proposalCells = @("F28:U28", "F20:U20", "H13:N13", "O7:U7", "F22:U22", "F24:U24", "F26:U26", "O13:U13")
foreach ($c in $proposalCells) {
    $proposal.Range($c).UnMerge()
    $proposal.Range($c).Merge()
}

# --- Reorder merged cells on "Flow Chart Template Tables" sheet ---
# Move the row-20/row-29/row-2/row-11 merge groups after the
# row-56/row-38/row-47 groups.
$flowChartCells = @("C20:F20", "G20:K20", "L20:O20", "C29:F29", "G29:J29", "K29:O29", "C2:F2", "G2:J2", "K2:N2", "C11:G11", "H11:K11", "L11:O11")
foreach ($c in $flowChartCells) {
    $flowChartTables.Range($c).UnMerge()
    $flowChartTables.Range($c).Merge()
}

# --- Move the selected/active tab from "Flow Chart" to "Flow Chart Template Tables" ---
$flowChartTables.Activate()
